# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Reorder "Santa Lucia" ahead of "Nueva Caledonia" in the country list
#    (row 207 was Nueva Caledonia, row 208 was Santa Lucia -> swap them)
$countryA = $ws.Range("A207").Value()
$countryB = $ws.Range("A208").Value()
$ws.Range("A207").Value = $countryB
$ws.Range("A208").Value = $countryA

# 2) Refresh the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 08:47"

# 3) Update the per-country stats that changed
#    columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
#             E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7447693
$ws.Range("C4").Value = 411
$ws.Range("D4").Value = 4700746
$ws.Range("E4").Value = 2535195
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 211752

# Row 27 - Israel
$ws.Range("B27").Value = 248133
$ws.Range("C27").Value = 2639
$ws.Range("D27").Value = 177752
$ws.Range("E27").Value = 68810
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 1571

# Row 28 - Ucrania
$ws.Range("B28").Value = 213028
$ws.Range("C28").Value = 4069
$ws.Range("D28").Value = 94443
$ws.Range("E28").Value = 114392
$ws.Range("G28").Value = 64
$ws.Range("H28").Value = 4193

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 56997
$ws.Range("C60").Value = 280
$ws.Range("D60").Value = 53457
$ws.Range("E60").Value = 3069
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 471

# Row 77 - El Salvador
$ws.Range("E77").Value = 4416
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = 848

# Row 80 - Australia
$ws.Range("B80").Value = 27096
$ws.Range("C80").Value = 18
$ws.Range("D80").Value = 24784
$ws.Range("E80").Value = 1424

# Row 175 - Taiwan
$ws.Range("B175").Value = 515
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 484
